# Product Offerings.docx — "Added business consulting to product offerings"
#
# After the final "Static Site" blurb ("This is for someone who wants a
# cheap site built quickly.") append a new "Business Consulting" section:
#   - Heading1 paragraph: "Business Consulting"
#   - Bold paragraph: "Expensive"
#   - Body paragraph describing the offering
#
# The trailing "_GoBack" bookmark (an empty/collapsed bookmark Word drops
# at the last edited spot) must end up inside the new final paragraph,
# positioned after its run — exactly where it sat relative to the old
# final paragraph before this edit. Re-creating a collapsed bookmark via
# Bookmarks.Add at a paragraph boundary in this host always serializes it
# *before* the adjoining run, so we instead insert the whole new block as
# literal OOXML (which preserves the exact element order we ask for) and
# delete/re-emit the bookmark as part of that literal markup.

$d = $word.ActiveDocument

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$docEnd = $d.Content.End - 1
$insertionPoint = $d.Range($docEnd, $docEnd)

$newBlockOoxml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Business Consulting</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Expensive</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Businesses and individual clients get business consultation from us, but could be redirected to David if they are in need of HUB services</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

[void]$insertionPoint.InsertXML($newBlockOoxml)
